$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.804.12"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "'2.279.98"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'504.26"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").Value = "'128.74"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'2.297.22"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "'0.0966"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "'0.344"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("D13").Value = "'4.94"
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").Value = "'23.35"
$ws.Range("E14").Value = "  +5.48%  "
$ws.Range("D15").Value = "'2.685.48"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "'54.849.37"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'0.0000131"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'2.292.06"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'10.31"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "'4.14"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'307.07"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "'6.45"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'60.03"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "'0.993"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'7.45"
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").Value = "'170.74"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'6.07"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'17.94"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "'0.995"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'3.79"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "'36.43"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'5.05"
$ws.Range("E42").Value = "  +6.32%  "
$ws.Range("D43").Value = "'3.40"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "'125.91"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'0.0498"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "'249.07"
$ws.Range("E46").Value = "  +4.88%  "
$ws.Range("D47").Value = "'0.0902"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "'0.549"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").Value = "'0.375"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "'10.82"
$ws.Range("E51").Value = "  +0.53%  "
